$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 holds a zero-padded code ("001"); force text format so Excel doesn't
# coerce it to the number 1 and drop the leading zeros.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"

$ws.Range("M2").Value = "2020-12-16 00:00:00"
$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 231658065.08
$ws.Range("P2").Value = 302.322625051
$ws.Range("Q2").Value = 1531175759.31
$ws.Range("R2").Value = 1998.2428619925
$ws.Range("S2").Value = 47870822.12
$ws.Range("T2").Value = 62.4732516939
$ws.Range("U2").Value = -89985223.72
$ws.Range("V2").Value = -117.4341546945

$ws.Range("Y2").Value = 44331270.07
$ws.Range("Z2").Value = 57.8540010458
$ws.Range("AA2").Value = -53879334.8
$ws.Range("AB2").Value = -70.3145902868
$ws.Range("AC2").Value = 76626109.29000001
$ws.Range("AD2").Value = 840.0640038649
